# F03 Froze Encoder 123
# Updates the per-epoch accuracy values in column B (rows 4-118, the rows
# whose values actually changed between the two training runs), refreshes
# the "DisplayOutputs" object-repr text in column A for rows 102-118 to the
# new process/object address, and restores the view state (selection /
# scroll anchor) to match where the sheet was left after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated accuracy values (row -> new value) -----------------------
$accuracyChanges = @{
    4   = 0.96875
    5   = 0.953125
    6   = 0.953125
    8   = 0.90625
    9   = 0.84375
    10  = 0.921875
    11  = 0.8125
    12  = 0.84375
    14  = 0.84375
    15  = 0.828125
    16  = 0.75
    18  = 0.796875
    19  = 0.8125
    20  = 0.796875
    21  = 0.78125
    22  = 0.78125
    23  = 0.78125
    24  = 0.765625
    25  = 0.78125
    26  = 0.78125
    27  = 0.796875
    28  = 0.796875
    29  = 0.796875
    32  = 0.8125
    34  = 0.8125
    43  = 0.796875
    44  = 0.796875
    45  = 0.796875
    46  = 0.796875
    48  = 0.78125
    61  = 0.78125
    62  = 0.78125
    63  = 0.78125
    64  = 0.78125
    65  = 0.78125
    66  = 0.78125
    67  = 0.78125
    68  = 0.78125
    69  = 0.78125
    70  = 0.78125
    71  = 0.78125
    72  = 0.78125
    73  = 0.78125
    103 = 0.828125
    104 = 0.8125
    105 = 0.796875
    106 = 0.796875
    107 = 0.765625
    108 = 0.765625
    109 = 0.75
    111 = 0.75
    112 = 0.78125
    113 = 0.703125
    114 = 0.6875
    116 = 0.78125
    117 = 0.90625
    118 = 0.7049180327868853
}

foreach ($row in $accuracyChanges.Keys) {
    $ws.Range("B$row").Value = $accuracyChanges[$row]
}

# --- Refresh the repr() text for the DisplayOutputs object in column A -
# Rows 102-118 all held the same repr string; the object's memory address
# changed on the new interpreter run.
$newRepr = "<__main__.DisplayOutputs object at 0x7f565065bb50>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Range("A$row").Value = $newRepr
}

# --- View state: select the data range used for the chart/analysis -----
$ws.Activate()
$ws.Range("A2:B118").Select()
